$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing SIM header from column B to column C, then
# put the new LINEID header into column B.
$ws.Range("C1").Value = $ws.Range("B1").Text
$ws.Range("B1").Value = "LINEID"

# Match the target column widths (B: 11 chars: C: ~22.43 chars, best-fit).
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws.Columns.Item(3).ColumnWidth = 21.666666666666668

# Move the selection to A2, matching the saved view state.
$ws.Range("A2").Select()

# Restore explicit portrait page setup.
$ws.PageSetup.Orientation = 1
